# Apply the "Commit ảnh cho các UC, Update các file doc đã làm" changes:
# - CODE-TC-SRS sheet: several Status (I) / Ghep giao dien (J) cells move to
#   "Done" (or one to "In Progress") along with their conditional-style fill.
# - DOC sheet: G12 moves to "Done"; F45 PIC changes from HaiCM to HuyenPT.
# - TONG HOP (summary) sheet counters are formulas; they recalculate
#   automatically once the underlying data above is updated.

$xlPasteFormats = -4122

$wb = $excel.ActiveWorkbook
$codeWs = $wb.Worksheets.Item("CODE-TC-SRS")
$docWs  = $wb.Worksheets.Item("DOC")

function Set-StatusCell {
    param($Worksheet, $TargetRef, $SourceRef, $Value)
    $src = $Worksheet.Range($SourceRef)
    $dst = $Worksheet.Range($TargetRef)
    $src.Copy() | Out-Null
    $dst.PasteSpecial($xlPasteFormats) | Out-Null
    $dst.Value2 = $Value
}

# --- CODE-TC-SRS sheet -------------------------------------------------
# "Done" cells - copy formatting from I7 (style index 4, the plain "Done" look)
Set-StatusCell $codeWs "I16" "I7" "Done"
Set-StatusCell $codeWs "I19" "I7" "Done"
Set-StatusCell $codeWs "I24" "I7" "Done"
Set-StatusCell $codeWs "J48" "I7" "Done"
Set-StatusCell $codeWs "I51" "I7" "Done"
Set-StatusCell $codeWs "I65" "I7" "Done"
Set-StatusCell $codeWs "I68" "I7" "Done"

# "Done" cells using the left-aligned "Done" look (style index 12), as in J7
Set-StatusCell $codeWs "J16" "J7" "Done"
Set-StatusCell $codeWs "J37" "J7" "Done"
Set-StatusCell $codeWs "J38" "J7" "Done"
Set-StatusCell $codeWs "J39" "J7" "Done"
Set-StatusCell $codeWs "J40" "J7" "Done"
Set-StatusCell $codeWs "J41" "J7" "Done"

# I70 -> Done using the centered "Done" look (style index 55), as in J70
Set-StatusCell $codeWs "I70" "J70" "Done"

# J59 -> "In Progress" using the bordered look (style index 17), as in I59
Set-StatusCell $codeWs "J59" "I59" "In Progress"

# --- DOC sheet -----------------------------------------------------------
# G12 -> "Done" (style index 4), as in G5
Set-StatusCell $docWs "G12" "G5" "Done"

# F45 PIC changed from HaiCM to HuyenPT (keep existing style/formatting)
$docWs.Range("F45").Value2 = "HuyenPT"

$excel.CutCopyMode = $false
